$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $NewValue, $ForceText)
    $cell = $Sheet.Range($CellRef)
    if ($ForceText) {
        $cell.Value = "'" + $NewValue
        $cell.Style = "Normal"
    } else {
        $cell.Value = $NewValue
    }
}

Set-TextValue $ws 'D2' '98.915.89' $false
Set-TextValue $ws 'E2' '  +0.68%  ' $false
Set-TextValue $ws 'D3' '3.292.16' $false
Set-TextValue $ws 'E3' '  -1.97%  ' $false
Set-TextValue $ws 'E4' '  +0.08%  ' $false
Set-TextValue $ws 'D5' '253.67' $true
Set-TextValue $ws 'E5' '  -1.96%  ' $false
Set-TextValue $ws 'D6' '623.89' $true
Set-TextValue $ws 'E6' '  +0.32%  ' $false
Set-TextValue $ws 'D7' '1.44' $true
Set-TextValue $ws 'E7' '  +19.45%  ' $false
Set-TextValue $ws 'E8' '  +6.21%  ' $false
Set-TextValue $ws 'E9' '  +0.01%  ' $false
Set-TextValue $ws 'D10' '0.977' $true
Set-TextValue $ws 'E10' '  +22.30%  ' $false
Set-TextValue $ws 'D11' '3.289.68' $false
Set-TextValue $ws 'E11' '  -1.96%  ' $false
Set-TextValue $ws 'D12' '0.200' $true
Set-TextValue $ws 'E12' '  +0.39%  ' $false
Set-TextValue $ws 'D13' '39.48' $true
Set-TextValue $ws 'E13' '  +10.14%  ' $false
Set-TextValue $ws 'D14' '98.540.92' $false
Set-TextValue $ws 'E14' '  +0.59%  ' $false
Set-TextValue $ws 'E15' '  +1.14%  ' $false
Set-TextValue $ws 'D16' '3.915.58' $false
Set-TextValue $ws 'E16' '  -1.19%  ' $false
Set-TextValue $ws 'D17' '5.48' $true
Set-TextValue $ws 'E17' '  -0.42%  ' $false
Set-TextValue $ws 'D18' '3.294.95' $false
Set-TextValue $ws 'E18' '  -1.83%  ' $false
Set-TextValue $ws 'D19' '3.47' $true
Set-TextValue $ws 'E19' '  -4.42%  ' $false
Set-TextValue $ws 'D20' '15.37' $true
Set-TextValue $ws 'E20' '  +2.51%  ' $false
Set-TextValue $ws 'D21' '6.33' $true
Set-TextValue $ws 'E21' '  +8.37%  ' $false
Set-TextValue $ws 'D22' '487.34' $true
Set-TextValue $ws 'E22' '  +0.95%  ' $false
Set-TextValue $ws 'D23' '9.34' $true
Set-TextValue $ws 'E23' '  +1.73%  ' $false
Set-TextValue $ws 'E24' '  -1.42%  ' $false
Set-TextValue $ws 'E25' '  -0.09%  ' $false
Set-TextValue $ws 'D26' '89.13' $true
Set-TextValue $ws 'E26' '  +0.97%  ' $false
Set-TextValue $ws 'D27' '0.326' $true
Set-TextValue $ws 'E27' '  +29.47%  ' $false
Set-TextValue $ws 'D28' '12.03' $true
Set-TextValue $ws 'E28' '  -0.39%  ' $false
Set-TextValue $ws 'D29' '3.472.55' $false
Set-TextValue $ws 'E29' '  -2.35%  ' $false
Set-TextValue $ws 'E30' '  -0.06%  ' $false
Set-TextValue $ws 'D31' '0.139' $true
Set-TextValue $ws 'E31' '  +14.98%  ' $false
Set-TextValue $ws 'E32' '  +3.16%  ' $false
Set-TextValue $ws 'D33' '10.33' $true
Set-TextValue $ws 'E33' '  +11.18%  ' $false
Set-TextValue $ws 'D34' '0.999' $true
Set-TextValue $ws 'E34' '  -0.02%  ' $false
Set-TextValue $ws 'D35' '27.91' $true
Set-TextValue $ws 'E35' '  +2.59%  ' $false
Set-TextValue $ws 'D36' '0.479' $true
Set-TextValue $ws 'E36' '  +7.09%  ' $false
Set-TextValue $ws 'E37' '  -1.38%  ' $false
Set-TextValue $ws 'E38' '  -2.36%  ' $false
Set-TextValue $ws 'E39' '  +0.33%  ' $false
Set-TextValue $ws 'D40' '24.79' $true
Set-TextValue $ws 'E40' '  -0.37%  ' $false
Set-TextValue $ws 'D41' '490.09' $true
Set-TextValue $ws 'E41' '  -4.78%  ' $false
Set-TextValue $ws 'D42' '3.63' $true
Set-TextValue $ws 'E42' '  +0.52%  ' $false
Set-TextValue $ws 'E43' '  -2.53%  ' $false
Set-TextValue $ws 'B44' 'USDe' $false
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/exbfr2U-0+usde-usde' $false
Set-TextValue $ws 'D44' '1.00' $true
Set-TextValue $ws 'E44' '  +0.01%  ' $false
Set-TextValue $ws 'B45' 'ARBITRUM' $false
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' $false
Set-TextValue $ws 'D45' '0.780' $true
Set-TextValue $ws 'E45' '  +0.39%  ' $false
Set-TextValue $ws 'E46' '  -4.87%  ' $false
Set-TextValue $ws 'B47' 'Stacks' $false
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' $false
Set-TextValue $ws 'D47' '1.95' $true
Set-TextValue $ws 'E47' '  +1.44%  ' $false
Set-TextValue $ws 'B48' 'Monero' $false
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' $false
Set-TextValue $ws 'D48' '157.79' $true
Set-TextValue $ws 'E48' '  -1.56%  ' $false
Set-TextValue $ws 'D49' '7.32' $true
Set-TextValue $ws 'E49' '  +16.15%  ' $false
Set-TextValue $ws 'B50' 'Filecoin' $false
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' $false
Set-TextValue $ws 'D50' '4.74' $true
Set-TextValue $ws 'E50' '  +5.00%  ' $false
Set-TextValue $ws 'B51' 'Mantle' $false
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' $false
Set-TextValue $ws 'D51' '0.843' $true
Set-TextValue $ws 'E51' '  +6.02%  ' $false
